$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target Price (D) value looks like a plain number must be
# forced to remain text (matching the source inlineStr cells), since Excel
# would otherwise auto-convert a numeric-looking string into a real number.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "41.877.08"
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").Value = "2.271.48"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue $ws.Range("D5") "304.19"
$ws.Range("E5").Value = "  +0.50%  "

Set-TextValue $ws.Range("D6") "92.93"
$ws.Range("E6").Value = "  +1.06%  "

$ws.Range("E7").Value = "  +1.98%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +0.28%  "

Set-TextValue $ws.Range("D10") "32.74"
$ws.Range("E10").Value = "  +2.10%  "

Set-TextValue $ws.Range("D11") "53.56"
$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("E13").Value = "  -1.22%  "

$ws.Range("E14").Value = "  +1.32%  "

$ws.Range("D15").Value = "2.624.06"
$ws.Range("E15").Value = "  +0.80%  "

$ws.Range("E16").Value = "  +0.79%  "

$ws.Range("D17").Value = "2.268.40"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("E18").Value = "  +3.50%  "

$ws.Range("D19").Value = "41.792.92"
$ws.Range("E19").Value = "  +1.26%  "

Set-TextValue $ws.Range("D20") "12.74"
$ws.Range("E20").Value = "  +4.41%  "

$ws.Range("E21").Value = "  +0.26%  "

Set-TextValue $ws.Range("D22") "5.94"
$ws.Range("E22").Value = "  +0.51%  "

Set-TextValue $ws.Range("D23") "67.16"
$ws.Range("E23").Value = "  +0.30%  "

Set-TextValue $ws.Range("D24") "243.62"
$ws.Range("E24").Value = "  +1.45%  "

$ws.Range("E25").Value = "  +0.32%  "

$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D26") "1.93"
$ws.Range("E26").Value = "  +3.71%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("E28").Value = "  +1.50%  "

Set-TextValue $ws.Range("D29") "9.54"
$ws.Range("E29").Value = "  -1.06%  "

Set-TextValue $ws.Range("D30") "2.08"
$ws.Range("E30").Value = "  -5.11%  "

Set-TextValue $ws.Range("D31") "35.32"
$ws.Range("E31").Value = "  +3.66%  "

Set-TextValue $ws.Range("D32") "160.94"
$ws.Range("E32").Value = "  +2.15%  "

$ws.Range("E33").Value = "  +1.38%  "

$ws.Range("E34").Value = "  -0.05%  "

Set-TextValue $ws.Range("D35") "0.0745"
$ws.Range("E35").Value = "  +1.07%  "

$ws.Range("E36").Value = "  -0.52%  "

Set-TextValue $ws.Range("D37") "17.12"
$ws.Range("E37").Value = "  +3.22%  "

$ws.Range("E38").Value = "  +2.65%  "

$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("E40").Value = "  +0.63%  "

$ws.Range("E41").Value = "  +1.73%  "

$ws.Range("E42").Value = "  -1.44%  "

$ws.Range("D43").Value = "2.006.70"
$ws.Range("E43").Value = "  -2.84%  "

Set-TextValue $ws.Range("D44") "19.57"
$ws.Range("E44").Value = "  -3.50%  "

$ws.Range("E45").Value = "  +2.00%  "

Set-TextValue $ws.Range("D46") "10.32"
$ws.Range("E46").Value = "  +2.13%  "

$ws.Range("E47").Value = "  +2.91%  "

$ws.Range("E48").Value = "  -1.91%  "

Set-TextValue $ws.Range("D49") "52.95"
$ws.Range("E49").Value = "  +3.56%  "

$ws.Range("E50").Value = "  +0.43%  "

$ws.Range("E51").Value = "  +1.06%  "
